# Append: 2025-10-18 12:41 JST
# Update the "取得日時" (acquisition timestamp) column (A) for the existing
# data rows (2-9) on the "ランサーズ" sheet from the previous run time
# (2025-10-18 12:33:02) to the new run time (2025-10-18 12:41:12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-18 12:33:02"
$newTimestamp = "2025-10-18 12:41:12"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
